# Weekly update: insert a new price record for Cebollín at
# Vega Monumental Concepción as row 72, pushing the existing
# rows 72-125 down to 73-126.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 72 (existing rows shift down).
$ws.Rows(72).Insert()

# Fill in the new record.
$ws.Range("A72").Value = 11
$ws.Range("B72").Value = "Vega Monumental Concepción"
$ws.Range("C72").Value = "Bíobío"
$ws.Range("D72").Value = 45126
$ws.Range("D72").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E72").Value = 8
$ws.Range("F72").Value = 100112037
$ws.Range("G72").Value = "Cebollín"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 80
$ws.Range("K72").Value = 4500
$ws.Range("L72").Value = 4500
$ws.Range("M72").Value = 4500
$ws.Range("N72").Value = "$/paquete 36 unidades"
$ws.Range("O72").Value = "Región Metropolitana"
$ws.Range("P72").Value = 125
$ws.Range("Q72").Value = 36
$ws.Range("R72").Value = "Hortaliza"
